$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert a new "Operation" column at the very left (column A).
#    Everything that used to live in A:AN now lives in B:AO.
# ------------------------------------------------------------------
$ws.Columns("A").Insert()

# ------------------------------------------------------------------
# 2. New header cell for the inserted column.
# ------------------------------------------------------------------
$ws.Range("A1").Value = "Operation"

# ------------------------------------------------------------------
# 3. Row 2 ("TEST - Dummy 15" product) is an Add operation - the rest
#    of the row already shifted into place correctly.
# ------------------------------------------------------------------
$ws.Range("A2").Value = "Add"

# ------------------------------------------------------------------
# 4. Row 3 used to be a full second product record; it is reworked to
#    be a "Change" operation that only carries the product identifier
#    plus the two fields that are actually changing (Turnaround Time
#    and Add to Inventory), with an explanatory note.
#    Clear everything out first, then lay down the new minimal data.
# ------------------------------------------------------------------
$ws.Range("B3:AO3").ClearContents()

$ws.Range("A3").Value = "Change"
$ws.Range("B3").Value = "Test - Dummy Product 2"
$ws.Range("C3").Value = "Dummy Product # 2"
$ws.Range("D3").Value = "TEST - Dummy 02"
$ws.Range("E3").Value = "In this test, change TA time to 7 days, and add 25 to inventory count."
$ws.Range("K3").Value = 7
$ws.Range("T3").Value = 25

# S3 keeps the hyperlink-style formatting (style used to carry the
# mailto hyperlink) but no longer has a value or an actual hyperlink.
$ws.Range("S3").Style = "Hyperlink"
$ws.Range("S3").ClearContents()

# ------------------------------------------------------------------
# 5. Remove the hyperlink that used to live on the old row 3 (it moved
#    with the data that got cleared above); keep the row 2 hyperlink.
# ------------------------------------------------------------------
foreach ($h in @($ws.Hyperlinks)) {
    if ($h.Range.Row -eq 3) {
        $h.Delete()
    }
}

# ------------------------------------------------------------------
# 6. Column widths: column A (Operation) gets a sensible width, and
#    column E (Brief Description) is widened to fit the longer note
#    text that can now appear there.
# ------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 9.1640625
$ws.Columns("E").ColumnWidth = 38.42

# ------------------------------------------------------------------
# 7. Selection moves to E4 after the edit.
# ------------------------------------------------------------------
$ws.Range("E4").Select()
